# Updated symbol list on Wed Jan 25 18:18:10 UTC 2023 with GitHub Actions
# Applies cell-level text updates (price/volume/hour refresh + a couple of
# coin re-rankings) to the crypto tracker worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '301.35'
    'E2' = '-3.02%'
    'G2' = '18'
    'D3' = '35.42'
    'E3' = '-0.19%'
    'G3' = '18'
    'D4' = '5.069'
    'E4' = '-0.86%'
    'G4' = '18'
    'D5' = '0.07988'
    'E5' = '-2.53%'
    'G5' = '18'
    'D6' = '1.883'
    'E6' = '-8.82%'
    'G6' = '18'
    'D7' = '7.787'
    'E7' = '-1.98%'
    'G7' = '18'
    'B8' = 'GateToken'
    'C8' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D8' = '4.048'
    'E8' = '-1.94%'
    'G8' = '18'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9279'
    'E9' = '0.28%'
    'G9' = '18'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.1439'
    'E10' = '28.65%'
    'G10' = '18'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1916'
    'E11' = '0.28%'
    'G11' = '18'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09058'
    'E12' = '-2.72%'
    'G12' = '18'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03500'
    'E13' = '-4.01%'
    'G13' = '18'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.09843'
    'E14' = '-0.58%'
    'G14' = '18'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001399'
    'E15' = '-2.17%'
    'G15' = '18'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005871'
    'E16' = '3.02%'
    'G16' = '18'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.533'
    'E17' = '1.77%'
    'G17' = '18'
    'D18' = '2.960'
    'E18' = '3.77%'
    'G18' = '18'
    'D19' = '0.3424'
    'E19' = '0.81%'
    'G19' = '18'
    'E20' = '0.17%'
    'G20' = '18'
    'D21' = '5.057'
    'E21' = '-0.60%'
    'G21' = '18'
    'D22' = '0.2398'
    'E22' = '8.56%'
    'G22' = '18'
    'D23' = '0.04494'
    'E23' = '-1.09%'
    'G23' = '18'
    'D24' = '0.001211'
    'E24' = '-1.07%'
    'G24' = '18'
    'D25' = '0.004764'
    'E25' = '-0.87%'
    'G25' = '18'
    'D26' = '0.0001230'
    'E26' = '-1.48%'
    'G26' = '18'
    'D27' = '0.0003024'
    'E27' = '-31.93%'
    'G27' = '18'
    'G28' = '18'
    'G29' = '18'
    'G30' = '18'
    'G31' = '18'
    'G32' = '18'
    'G33' = '18'
    'G34' = '18'
    'G35' = '18'
    'G36' = '18'
    'G37' = '18'
    'G38' = '18'
    'D39' = '0.01833'
    'E39' = '-7.06%'
    'G39' = '18'
    'D40' = '0.04764'
    'E40' = '-2.75%'
    'G40' = '18'
    'D41' = '0.01054'
    'E41' = '14.15%'
    'G41' = '18'
    'E42' = '-3.28%'
    'G42' = '18'
    'D43' = '0.1326'
    'E43' = '-4.12%'
    'G43' = '18'
    'D44' = '0.002110'
    'E44' = '-3.80%'
    'G44' = '18'
    'D45' = '0.01100'
    'E45' = '-5.32%'
    'G45' = '18'
    'D46' = '0.00006229'
    'E46' = '-4.82%'
    'G46' = '18'
    'D47' = '0.00000000750'
    'E47' = '0.13%'
    'G47' = '18'
    'D48' = '64.66'
    'E48' = '-64.07%'
    'G48' = '18'
    'E49' = '10.78%'
    'G49' = '18'
    'D50' = '0.00002100'
    'E50' = '0.13%'
    'G50' = '18'
    'D51' = '0.0002000'
    'E51' = '0.13%'
    'G51' = '18'
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    # Force text storage so numeric-looking strings (prices, percentages,
    # the hour value) are not silently converted to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$key]
}
